$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "63.575.53"
$ws.Cells.Item(2, 5).Value = "  +3.32%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.070.94"
$ws.Cells.Item(3, 5).Value = "  +2.29%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "551.26"
$ws.Cells.Item(5, 5).Value = "  +2.45%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "142.48"
$ws.Cells.Item(6, 5).Value = "  +5.35%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.999"
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "3.069.28"
$ws.Cells.Item(8, 5).Value = "  +2.37%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.503"
$ws.Cells.Item(9, 5).Value = "  +1.30%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "6.55"
$ws.Cells.Item(10, 5).Value = "  +6.45%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +2.51%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.458"
$ws.Cells.Item(12, 5).Value = "  +2.15%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000229"
$ws.Cells.Item(13, 5).Value = "  +2.69%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "35.04"
$ws.Cells.Item(14, 5).Value = "  +2.80%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.560.77"
$ws.Cells.Item(15, 5).Value = "  +2.21%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "63.480.80"
$ws.Cells.Item(16, 5).Value = "  +3.21%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.065.46"
$ws.Cells.Item(17, 5).Value = "  +2.45%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -0.92%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.81"
$ws.Cells.Item(19, 5).Value = "  +2.37%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "487.06"
$ws.Cells.Item(20, 5).Value = "  +4.14%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.92"
$ws.Cells.Item(21, 5).Value = "  +5.01%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.680"
$ws.Cells.Item(22, 5).Value = "  +0.33%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "7.30"
$ws.Cells.Item(23, 5).Value = "  +4.78%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "81.28"
$ws.Cells.Item(24, 5).Value = "  +1.17%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "12.77"
$ws.Cells.Item(25, 5).Value = "  +6.33%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.00"
$ws.Cells.Item(26, 5).Value = "  -0.18%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.77"
$ws.Cells.Item(27, 5).Value = "  +3.04%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.93"
$ws.Cells.Item(28, 5).Value = "  +1.59%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.02"
$ws.Cells.Item(29, 5).Value = "  +7.03%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.997"
$ws.Cells.Item(30, 5).Value = "  -0.04%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "26.32"
$ws.Cells.Item(31, 5).Value = "  +2.58%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.17"
$ws.Cells.Item(32, 5).Value = "  +1.29%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.46"
$ws.Cells.Item(33, 5).Value = "  +7.87%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.72"
$ws.Cells.Item(34, 5).Value = "  +3.15%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "55.56"
$ws.Cells.Item(35, 5).Value = "  +1.42%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.02"
$ws.Cells.Item(36, 5).Value = "  +1.83%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "467.03"
$ws.Cells.Item(37, 5).Value = "  +2.87%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0827"
$ws.Cells.Item(38, 5).Value = "  +4.71%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0400"
$ws.Cells.Item(39, 5).Value = "  +3.56%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.048.16"
$ws.Cells.Item(40, 5).Value = "  -3.88%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -1.38%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "8.26"
$ws.Cells.Item(42, 5).Value = "  +1.30%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.59"
$ws.Cells.Item(43, 5).Value = "  +4.57%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "27.93"
$ws.Cells.Item(44, 5).Value = "  +2.96%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.257"
$ws.Cells.Item(45, 5).Value = "  +4.68%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.09%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.06"
$ws.Cells.Item(47, 5).Value = "  +2.94%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.111"
$ws.Cells.Item(48, 5).Value = "  +2.47%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "PEPE"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0₃0513"
$ws.Cells.Item(49, 5).Value = "  +2.92%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Monero"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "116.73"
$ws.Cells.Item(50, 5).Value = "  -1.73%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.09"
$ws.Cells.Item(51, 5).Value = "  +4.07%  "
